$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Ver video tutorial de GitHub."

$ws.Range("AMH6:AMJ9").Clear()

$ws.Range("B2").Select()
